# Update the latest cryptocurrency price/volume snapshot values
# (columns D = Price, E = Volume(1h)) on the active worksheet.
# Values are written with a leading single-quote where needed so that
# Excel keeps them as text (matching the sheet's original inlineStr
# formatting) instead of auto-converting them to numbers, which would
# silently drop significant trailing/leading zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.390.17'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.848.77'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''240.44'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '''0.6276'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '''0.07632'
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('D9').Value = '''0.2906'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('D10').Value = '''24.72'
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').Value = '''0.07739'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '''5.020'
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '''0.6792'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').Value = '''0.00001061'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '''83.19'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '''6.158'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '29.405.92'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '''226.80'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '''12.33'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '''7.497'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').Value = '''1.0000'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '''158.10'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').Value = '''0.1379'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').Value = '''8.403'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').Value = '''1.387'
$ws.Range('E27').Value = '  +5.58%  '
$ws.Range('D28').Value = '''1.464'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').Value = '''0.05598'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''4.121'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').Value = '''4.075'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('D32').Value = '''1.837'
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D34').Value = '''0.6949'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('D35').Value = '''2.583'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').Value = '1.231.89'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').Value = '''0.01800'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Value = '''2.719'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').Value = '''6.400'
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('D40').Value = '''0.9055'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = '''1.001'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '''101.61'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '''66.05'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '''7.176'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '''0.00000000119'
$ws.Range('E45').Value = '  -2.80%  '
$ws.Range('D46').Value = '''0.4015'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = '''8.987'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '''1.681'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  +1.93%  '
$ws.Range('D50').Value = '''0.05703'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('E51').Value = '  +0.06%  '
